$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fix row 61 height: remove the stray custom row height (ht="30") ---
$ws.Rows.Item(61).AutoFit()

# --- Grow Table1 by two new data rows plus one trailing blank row ---
$lo = $ws.ListObjects.Item(1)
$lo.ListRows.Add() | Out-Null
$lo.ListRows.Add() | Out-Null
$lo.ListRows.Add() | Out-Null

# Put values in first so the paste-format step below doesn't clobber them
$ws.Range("A149").Value = "ايجاد، حذف، نمايش و ويرايش بخش zone هاي دسترسي كاربران"
$ws.Range("B149").Value = "دوم"
$ws.Range("C149").Value = 0
$ws.Range("D149").Value = 0
$ws.Range("E149").Value = 0
$ws.Range("F149").Value = 0

$ws.Range("A150").Value = "ايجاد، حذف، نمايش و ويرايش بخش componentهاي نرم افزار"
$ws.Range("B150").Value = "دوم"
$ws.Range("C150").Value = 0
$ws.Range("D150").Value = 0
$ws.Range("E150").Value = 0
$ws.Range("F150").Value = 0

# Copy the formatting of the existing last data row (148) onto the two
# new rows, and the formatting of a blank formatted row (125) onto the
# new trailing blank row (151) -- keeps the same look (wrap/centering)
# without bloating the style table.
$ws.Range("A148:F148").Copy()
$ws.Range("A149:F149").PasteSpecial(-4122)
$ws.Range("A148:F148").Copy()
$ws.Range("A150:F150").PasteSpecial(-4122)
$ws.Range("B125:F125").Copy()
$ws.Range("B151:F151").PasteSpecial(-4122)

# --- Update the active selection to reflect where the user ended up ---
$ws.Range("A151").Select() | Out-Null
